$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.488.43'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '3.394.62'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.65'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.83'
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.475'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.69'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.93%  '
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.388'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("D12").Value = '3.976.12'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.410.76'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = '61.496.27'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.14'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.68'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.54'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.42'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.555'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000113'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.193'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.38%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.04'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -3.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.43'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.94'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.38'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("D37").Value = '3.428.26'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0769'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.20'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.779'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.67'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = '2.458.58'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.09'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.71'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0263'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.07'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("E51").Value = '  -1.27%  '
